# Update de Plazo Fijo 12-01-2023
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Modulos" sheet: add the new "Plazo Fijo" menu entries below the
# existing data (rows 5-13, columns A-D), written in the same order the
# original author typed them in (this controls the order new shared
# strings are appended in).
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Modulos")

$ws.Range("A5").Value  = "Plazo Fijo"
$ws.Range("D5").Value  = "Liquidacion Plazo Fijo Ajustable"
$ws.Range("B5").Value  = "Impresión de Certificados y Formularios"
$ws.Range("B6").Value  = "Consulta de Plazo Fijo"
$ws.Range("D6").Value  = "Consulta Circulares Vigentes"
$ws.Range("D7").Value  = "Consulta Tasa Plazo Fijo"
$ws.Range("D8").Value  = "Consulta Indices"
$ws.Range("D9").Value  = "Altas/Pagos/Inm Pagados/Imm Impagos"
$ws.Range("D11").Value = "Plazo Fijo Inmovilizado"
$ws.Range("D10").Value = "PF Vencidos Inmovilizados Pagados"
$ws.Range("C12").Value = "Consulta por Firmante"
$ws.Range("D12").Value = "Monto Activo por Titular"
$ws.Range("D13").Value = "Plazo Fijos Activos"

# Update the view: scroll so row 4 is at the top and select D13, the
# last cell entered.
$ws.Activate()
$ws.Range("D13").Select()
$excel.ActiveWindow.ScrollRow = 4

# ---------------------------------------------------------------------
# "Users" sheet: selection moved from E9 to A2.
# ---------------------------------------------------------------------
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Activate()
$wsUsers.Range("A2").Select()

# Leave focus back on the Modulos sheet, which is the active tab.
$ws.Activate()
